# ITO-000 update API create User
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("endpoint")

# Insert a new row at 7, pushing current rows 7..10 down to 8..11
$ws.Rows.Item(7).Insert() | Out-Null

# Row 6: add Data-RequestBody sample JSON in H6, change response model in L6
$ws.Range("H6").Value = "{`n  ""email"":""hoangnhuocquy@csc.com"",`n  ""password"":""P@ssword123"",`n  ""displayName"":""Hoang Nhuoc QUy""`n}"
$ws.Range("H6").WrapText = $true
$ws.Range("H6").VerticalAlignment = -4160

$ws.Range("L6").Value = "onlinejudge.dto.MyResponse"

# Row 7: new "user already exists" response row
$ws.Range("H7").WrapText = $true
$ws.Range("H7").VerticalAlignment = -4160

$ws.Range("K7").Value = 400
$ws.Range("L7").Value = "onlinejudge.dto.MyResponse"
$ws.Range("N7").Value = "User exist"
$ws.Range("O7").Value = "user.create.exist"
$ws.Range("P7").Value = "User exist with email [{0}]"

# Row heights for rows 6 & 7 (taller to fit wrapped JSON sample)
$ws.Rows.Item(6).RowHeight = 26.25
$ws.Rows.Item(7).RowHeight = 26.25

# Update selection / view (scroll so column H is leftmost, select P7)
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 8
$ws.Range("P7").Select() | Out-Null

Write-Host "done"
